# Commit: "updated common dataprovider for valid,invalid,titlecheck"
#
# 1. invalidCredentialTest: trim the duplicated admin124..admin132 rows,
#    keeping just the Peter/Balaji negative-credential rows.
# 2. Insert a brand-new "validateTitleTest" sheet between validCredentialTest
#    and addEmployeeTest, holding the Title/OrangeHRM expected-title pair.
# 3. Leave view state (active tab/selection) the way Excel would after that
#    edit session: invalidCredentialTest ends up the active/selected tab.

$wb = $excel.ActiveWorkbook

# --- 1. invalidCredentialTest: drop rows 4:13 -------------------------------
$ws1 = $wb.Worksheets.Item("invalidCredentialTest")
$ws1.Activate()
$ws1.Rows("4:13").Delete()
[void]$ws1.Rows("4:1048576").Select()

# --- 2. validCredentialTest: visit it (matches the view-state left in the
#        sheet after the edit session) ---------------------------------------
$ws2 = $wb.Worksheets.Item("validCredentialTest")
$ws2.Activate()
[void]$ws2.Range("A3").EntireRow.Select()

# --- 3. New sheet "validateTitleTest" right after validCredentialTest -------
$newSheet = $wb.Worksheets.Add($null, $ws2)
$newSheet.Name = "validateTitleTest"
$newSheet.Range("A1").Value = "Title"
$newSheet.Range("A2").Value = "OrangeHRM"
$newSheet.Columns("A:A").AutoFit()
[void]$newSheet.Range("A3").Select()

# --- 4. Final active tab is invalidCredentialTest ---------------------------
$ws1.Activate()
